# Update AVTA Course Excel - add offshore Region1 Horticulture courses
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2 - AHC30716 Certificate III in Horticulture (entered in full,
# department filled in later together with the other rows)
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "AHC30716"
$ws.Range("B2").Value = "110597F"
$ws.Range("D2").Value = "CERTIFICATE III IN HORTICULTURE"
$ws.Range("E2").Value = 52
$ws.Range("H2").Value = "44 wks Tuition + 8 wks Break"
$ws.Range("I2").Value = 12200
$ws.Range("J2").Value = "12,000 tuition fee + 200 handling fee"
$ws.Range("M2").Value = "TAS"

# ---------------------------------------------------------------------------
# Row 3 - AHC40416 Certificate IV in Horticulture
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "AHC40416"
$ws.Range("B3").Value = "110598E"
$ws.Range("D3").Value = "CERTIFICATE IV IN HORTICULTURE"
$ws.Range("E3").Value = 52
$ws.Range("H3").Value = "44 wks Tuition + 8 wks Break"
$ws.Range("I3").Value = 14200
$ws.Range("M3").Value = "TAS"

# ---------------------------------------------------------------------------
# Row 4 - AHC51422 Diploma of Agribusiness Management
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "AHC51422"
$ws.Range("B4").Value = "110774E"
$ws.Range("D4").Value = "DIPLOMA OF AGRIBUSINESS MANAGEMENT"
$ws.Range("E4").Value = 52
$ws.Range("H4").Value = "44 wks Tuition + 8 wks Break"
$ws.Range("I4").Value = 15200
$ws.Range("M4").Value = "TAS"

# ---------------------------------------------------------------------------
# Row 5 - Package: Certificate III + Certificate IV in Horticulture
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "AHC30716 / AHC40416"
$ws.Range("B5").Value = "110597F / 110598E"
$ws.Range("D5").Value = "CERTIFICATE III IN HORTICULTURE +`nCERTIFICATE IV IN HORTICULTURE"
$ws.Range("E5").Value = 104
$ws.Range("H5").Value = "88 wks Tuition + 16 wks Break"
$ws.Range("I5").Value = 22200
$ws.Range("M5").Value = "TAS"

# ---------------------------------------------------------------------------
# Row 6 - Package: Certificate III Horticulture + Diploma Agribusiness Mgmt
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "AHC30716 / AHC51422"
$ws.Range("B6").Value = "110597F / 110774E"
$ws.Range("D6").Value = "CERTIFICATE III IN HORTICULTURE +`nDIPLOMA OF AGRIBUSINESS MANAGEMENT"
$ws.Range("E6").Value = 104
$ws.Range("H6").Value = "88 wks Tuition + 16 wks Break"
$ws.Range("I6").Value = 23200
$ws.Range("M6").Value = "TAS"

# ---------------------------------------------------------------------------
# Row 7 - Package: Certificate IV Horticulture + Diploma Agribusiness Mgmt
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "AHC40416 / AHC51422"
$ws.Range("B7").Value = "110598E / 110774E"
$ws.Range("D7").Value = "CERTIFICATE IV IN HORTICULTURE +`nDIPLOMA OF AGRIBUSINESS MANAGEMENT"
$ws.Range("E7").Value = 104
$ws.Range("H7").Value = "88 wks Tuition + 16 wks Break"
$ws.Range("I7").Value = 23200
$ws.Range("M7").Value = "TAS"

# ---------------------------------------------------------------------------
# Department column (C) - filled in for all rows after the rest of the
# table was populated
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "HORTICULTURE"
$ws.Range("C3").Value = "HORTICULTURE"
$ws.Range("C4").Value = "MANAGEMENT"
$ws.Range("C5").Value = "PACKAGES"
$ws.Range("C6").Value = "PACKAGES"
$ws.Range("C7").Value = "PACKAGES"

# ---------------------------------------------------------------------------
# Tuition detail column (J) - filled in for the remaining rows last
# ---------------------------------------------------------------------------
$ws.Range("J3").Value = "14,000 tuition fee + 200 handling fee"
$ws.Range("J4").Value = "15,000 tuition fee + 200 handling fee"
$ws.Range("J5").Value = "22,000 tuition fee + 200 handling fee"
$ws.Range("J6").Value = "23,000 tuition fee + 200 handling fee"
$ws.Range("J7").Value = "23,000 tuition fee + 200 handling fee"

# ---------------------------------------------------------------------------
# Row heights (45pt for all new rows, to match wrapped text content)
# ---------------------------------------------------------------------------
$ws.Range("A2:R7").RowHeight = 45

# ---------------------------------------------------------------------------
# Number formats - thousands separator on tuition amount column (I), and
# also applied (together with wrap) to the tuition detail column (J)
# ---------------------------------------------------------------------------
$ws.Range("I2:I7").NumberFormat = "#,##0"
$ws.Range("J2:J7").NumberFormat = "#,##0"

# ---------------------------------------------------------------------------
# Wrap text - column H (duration detail) and J (tuition detail) on all
# new rows, plus A, B, D on the package rows (5-7)
# ---------------------------------------------------------------------------
$ws.Range("H2:H7").WrapText = $true
$ws.Range("J2:J7").WrapText = $true
$ws.Range("A5:B7").WrapText = $true
$ws.Range("D5:D7").WrapText = $true

# ---------------------------------------------------------------------------
# Selection / view state
# ---------------------------------------------------------------------------
$null = $ws.Range("J15").Select()
